# Minor improvements in board routing. Generated CAM files. Fixed footprints in BoM.
#
# The only functional change is in the Bill of Materials worksheet: two
# footprint values that were erroneously "0604" / "0605" are corrected to
# the standard "0603" footprint used throughout the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 -> C43 footprint: 0604 -> 0603
$ws.Range("C3").Value = "0603"

# Row 4 -> C41 footprint: 0605 -> 0603
$ws.Range("C4").Value = "0603"

# Leave the cursor/selection where the author last clicked before saving.
$ws.Range("C10").Select()

$wb.Save()
